$wb = $excel.ActiveWorkbook

# --- start_price sheet ---
$wsStartPrice = $wb.Worksheets.Item("start_price")
$wsStartPrice.Range("A2").Value = 0.0151

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.000004896473920455704
$wsLinear.Range("B3").Value = -0.524192538636706
$wsLinear.Range("B4").Value = 0.00006065929019144551
$wsLinear.Range("B5").Value = "[0.9999999999999999, 0.3997991326778741, 0.1429103721719559, 0.0974873311610669, 0.08621518452083317, 0.16070029683432624, 0.26167602599698364, 0.2258810318067539, 0.16812639347177927, 0.10719200812474186, 0.06627767242589713, 0.0951537806060099, 0.1355668304722252, 0.17375006650385535, 0.24847442936540187, 0.1930895706587592, 0.12456878908065268, 0.11576667881684702, 0.05977547418034912, 0.06984261773849555]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 1.052580331061344
$wsNonLinear.Range("B4").Value = -0.0005379101107449344
$wsNonLinear.Range("B5").Value = -0.6261225404145367
$wsNonLinear.Range("B6").Value = 0.00005323110591756335
$wsNonLinear.Range("B7").Value = 0.000748671109656337
$wsNonLinear.Range("B8").Value = -0.7262946049081308
$wsNonLinear.Range("B9").Value = 0.0000680687399581896
$wsNonLinear.Range("B10").Value = "[1.0, 0.4019127238948999, 0.14587175322265739, 0.09855777458335774, 0.08740943477965374, 0.16647003976941224, 0.2685426803566058, 0.22670836394892155, 0.16746595313670706, 0.10565915107925201, 0.07098596179930973, 0.09702325790797552, 0.1380695973666252, 0.17440484664114184, 0.248204836449194, 0.20003269406855798, 0.1301715032398129, 0.11933110191376471, 0.06411004567701609, 0.06968750467004141]"
